# Updated symbol list on Mon Dec 26 17:00:05 UTC 2022 with GitHub Actions
#
# All data cells on this sheet were authored as literal text (inline
# strings), including the "Price" column which looks numeric but must
# keep its exact textual formatting (trailing zeros, precision, etc).
# Plain `.Value = "<numeric-looking-string>"` lets Excel's automatic
# type detection coerce it into a real number (losing e.g. "242.60" ->
# 242.6), so for every numeric-looking cell we first force the cell to
# Text format, then assign the string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Simple price refreshes (no row movement) ---
Set-TextValue "D2" "242.60"
Set-TextValue "D3" "23.02"
Set-TextValue "D4" "5.424"
Set-TextValue "D6" "3.438"
Set-TextValue "D7" "6.521"
Set-TextValue "D8" "0.8104"
Set-TextValue "D9" "0.9776"

# --- Rows 10-18: the coin that used to be at the bottom of this block
#     ("One") moved to the top, shifting WazirX..CoinExToken down by one
#     row each; prices/volume labels updated to the new snapshot ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01129"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1419"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07424"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03271"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03055"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09346"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.847"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001586"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04671"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- Remaining standalone price refreshes further down the sheet ---
Set-TextValue "D19" "0.005885"
Set-TextValue "D20" "0.001268"
Set-TextValue "D23" "3.592"
Set-TextValue "D27" "0.0002285"
Set-TextValue "D40" "0.03936"
Set-TextValue "D41" "0.006186"
Set-TextValue "D42" "0.1069"
Set-TextValue "D43" "0.003001"
Set-TextValue "D44" "0.009135"
$ws.Range("E44").Value = "43LocalTradersLCT"
Set-TextValue "D45" "0.00005198"
Set-TextValue "D47" "0.6661"
Set-TextValue "D48" "0.002384"
